$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 2549.8
$ws.Range("I13").Value = 916.3333
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 916.3333
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = -747.3333
$ws.Range("N13").Value = -5338

$ws.Range("H31").Value = 831.2857
$ws.Range("I31").Value = 569.8333
$ws.Range("J31").Value = 2400
$ws.Range("K31").Value = 1709.4999
$ws.Range("L31").Value = 7200
$ws.Range("M31").Value = -1479.4999
$ws.Range("N31").Value = -7660

$ws.Range("H43").Value = 4396589.5
$ws.Range("I43").Value = 17553058
$ws.Range("J43").Value = 11100.223
$ws.Range("K43").Value = 17553058
$ws.Range("L43").Value = 11100.223
$ws.Range("M43").Value = -17552989
$ws.Range("N43").Value = -11238.223

$ws.Range("H135").Value = 4897.423
$ws.Range("I135").Value = 1697.4445
$ws.Range("J135").Value = 12097.375
$ws.Range("K135").Value = 15277.0005
$ws.Range("L135").Value = 108876.375
$ws.Range("M135").Value = -12742.0005

$ws.Range("H136").Value = 76750
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 76750
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 76750
$ws.Range("N136").Value = -86950

$ws.Range("H137").Value = 1541.2084
$ws.Range("I137").Value = 1077.7646
$ws.Range("J137").Value = 2666.7144
$ws.Range("K137").Value = 3233.2938
$ws.Range("L137").Value = 8000.1432
$ws.Range("M137").Value = -683.2937999999999
$ws.Range("N137").Value = -13100.1432

$ws.Range("H138").Value = 1673833.4
$ws.Range("I138").Value = 4000
$ws.Range("J138").Value = 2508750
$ws.Range("K138").Value = 12000
$ws.Range("L138").Value = 7526250
$ws.Range("M138").Value = -6860
$ws.Range("N138").Value = -7536530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1839.6
$ws.Range("I2").Value = 1488.4445
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 1488.4445
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -1375.4445

$ws.Range("H45").Value = 4087.3333
$ws.Range("I45").Value = 2540.8572
$ws.Range("J45").Value = 9500
$ws.Range("K45").Value = 2540.8572
$ws.Range("L45").Value = 9500
$ws.Range("M45").Value = -2163.8572

$ws.Range("H74").Value = 1271.1666
$ws.Range("I74").Value = 1129.3334
$ws.Range("J74").Value = 1696.6666
$ws.Range("K74").Value = 1129.3334
$ws.Range("L74").Value = 1696.6666
$ws.Range("M74").Value = -255.3334
$ws.Range("N74").Value = -3444.6666

$ws.Range("H77").Value = 1271.1666
$ws.Range("I77").Value = 1129.3334
$ws.Range("J77").Value = 1696.6666
$ws.Range("K77").Value = 5646.666999999999
$ws.Range("L77").Value = 8483.333000000001
$ws.Range("M77").Value = -1278.666999999999
$ws.Range("N77").Value = -17219.333

$ws.Range("H116").Value = 1839.6
$ws.Range("I116").Value = 1488.4445
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 1488.4445
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 805.5554999999999

$ws.Range("H132").Value = 2960.48
$ws.Range("I132").Value = 1534.359
$ws.Range("J132").Value = 8016.727
$ws.Range("K132").Value = 4603.076999999999
$ws.Range("L132").Value = 24050.181
$ws.Range("M132").Value = -2073.076999999999
$ws.Range("N132").Value = -29110.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1839.6
$ws.Range("I3").Value = 1488.4445
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 1488.4445
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -1374.4445

$ws.Range("H95").Value = 60076.625
$ws.Range("I95").Value = 150000
$ws.Range("J95").Value = 47230.43
$ws.Range("K95").Value = 150000
$ws.Range("L95").Value = 47230.43
$ws.Range("M95").Value = -147254
$ws.Range("N95").Value = -52722.43

$ws.Range("H99").Value = 2299.6667
$ws.Range("I99").Value = 2299.6667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2299.6667
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -801.6667000000002

$ws.Range("H105").Value = 1280
$ws.Range("I105").Value = 800
$ws.Range("J105").Value = 3200
$ws.Range("K105").Value = 800
$ws.Range("L105").Value = 3200
$ws.Range("M105").Value = 947

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 242.77142
$ws.Range("I22").Value = 259.33334
$ws.Range("J22").Value = 217.92857
$ws.Range("K22").Value = 259.33334
$ws.Range("L22").Value = 217.92857
$ws.Range("M22").Value = 90.66665999999998
$ws.Range("N22").Value = -917.92857

$ws.Range("H31").Value = 1479.6123
$ws.Range("I31").Value = 1436.1915
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1436.1915
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -1141.1915
$ws.Range("N31").Value = -3090

$ws.Range("H34").Value = 1479.6123
$ws.Range("I34").Value = 1436.1915
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 1436.1915
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -1234.1915
$ws.Range("N34").Value = -2904

$ws.Range("H58").Value = 3256.6924
$ws.Range("I58").Value = 2165.2222
$ws.Range("J58").Value = 5712.5
$ws.Range("K58").Value = 2165.2222
$ws.Range("L58").Value = 5712.5
$ws.Range("M58").Value = -1962.2222

$ws.Range("H103").Value = 45005.5
$ws.Range("I103").Value = 40006.145
$ws.Range("J103").Value = 80001
$ws.Range("K103").Value = 40006.145
$ws.Range("L103").Value = 80001
$ws.Range("M103").Value = -38834.145

$ws.Range("H136").Value = 3256.6924
$ws.Range("I136").Value = 2165.2222
$ws.Range("J136").Value = 5712.5
$ws.Range("K136").Value = 6495.6666
$ws.Range("L136").Value = 17137.5
$ws.Range("M136").Value = -3945.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8071.5
$ws.Range("I5").Value = 470.375
$ws.Range("J5").Value = 18206.334
$ws.Range("K5").Value = 1411.125
$ws.Range("L5").Value = 54619.00199999999
$ws.Range("M5").Value = -1299.125
$ws.Range("N5").Value = -54843.00199999999

$ws.Range("H12").Value = 1308
$ws.Range("I12").Value = 2400
$ws.Range("J12").Value = 840
$ws.Range("K12").Value = 7200
$ws.Range("L12").Value = 2520
$ws.Range("M12").Value = -7027
$ws.Range("N12").Value = -2866

$ws.Range("H36").Value = 1795.3334
$ws.Range("I36").Value = 93
$ws.Range("J36").Value = 3923.25
$ws.Range("K36").Value = 279
$ws.Range("L36").Value = 11769.75
$ws.Range("M36").Value = -12107.75

$ws.Range("H113").Value = 591.8461
$ws.Range("I113").Value = 482.66666
$ws.Range("J113").Value = 837.5
$ws.Range("K113").Value = 1447.99998
$ws.Range("L113").Value = 2512.5
$ws.Range("M113").Value = 722.0000199999999
$ws.Range("N113").Value = -6852.5

$ws.Range("H122").Value = 946.5333000000001
$ws.Range("I122").Value = 789.5
$ws.Range("J122").Value = 970.6923
$ws.Range("K122").Value = 7105.5
$ws.Range("L122").Value = 8736.2307
$ws.Range("M122").Value = -4655.5
$ws.Range("N122").Value = -13636.2307

$ws.Range("H132").Value = 2216.7407
$ws.Range("I132").Value = 1240.5
$ws.Range("J132").Value = 2495.6667
$ws.Range("K132").Value = 11164.5
$ws.Range("L132").Value = 22461.0003
$ws.Range("M132").Value = -8634.5
$ws.Range("N132").Value = -27521.0003

$ws.Range("H135").Value = 8071.5
$ws.Range("I135").Value = 470.375
$ws.Range("J135").Value = 18206.334
$ws.Range("K135").Value = 4233.375
$ws.Range("L135").Value = 163857.006
$ws.Range("M135").Value = -1698.375
$ws.Range("N135").Value = -168927.006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 41622.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 41622.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 41622.5
$ws.Range("N68").Value = -43244.5

$ws.Range("H71").Value = 41622.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 41622.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 124867.5
$ws.Range("N71").Value = -132979.5

$ws.Range("H132").Value = 4869.881
$ws.Range("I132").Value = 4792.4478
$ws.Range("J132").Value = 5175.0586
$ws.Range("K132").Value = 14377.3434
$ws.Range("L132").Value = 15525.1758
$ws.Range("M132").Value = -11847.3434
$ws.Range("N132").Value = -20585.1758

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 13316.968
$ws.Range("I61").Value = 11002.28
$ws.Range("J61").Value = 22961.5
$ws.Range("K61").Value = 11002.28
$ws.Range("L61").Value = 22961.5
$ws.Range("M61").Value = -10800.28
$ws.Range("N61").Value = -23365.5

$ws.Range("H93").Value = 1721.138
$ws.Range("I93").Value = 1496.52
$ws.Range("J93").Value = 3125
$ws.Range("K93").Value = 1496.52
$ws.Range("L93").Value = 3125
$ws.Range("M93").Value = -248.52
$ws.Range("N93").Value = -5621

$ws.Range("H100").Value = 3658
$ws.Range("I100").Value = 3054.3
$ws.Range("J100").Value = 4520.4287
$ws.Range("K100").Value = 3054.3
$ws.Range("L100").Value = 4520.4287
$ws.Range("M100").Value = -2513.3
$ws.Range("N100").Value = -5602.4287

$ws.Range("H113").Value = 13316.968
$ws.Range("I113").Value = 11002.28
$ws.Range("J113").Value = 22961.5
$ws.Range("K113").Value = 11002.28
$ws.Range("L113").Value = 22961.5
$ws.Range("M113").Value = -8832.280000000001
$ws.Range("N113").Value = -27301.5

$ws.Range("H132").Value = 5391.8335
$ws.Range("I132").Value = 4305.7
$ws.Range("J132").Value = 6749.5
$ws.Range("K132").Value = 12917.1
$ws.Range("L132").Value = 20248.5
$ws.Range("M132").Value = -10387.1

$ws.Range("H134").Value = 103476.336
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 103476.336
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 103476.336
$ws.Range("N134").Value = -113616.336

$ws.Range("H136").Value = 4519.641
$ws.Range("I136").Value = 4553.5757
$ws.Range("J136").Value = 4333
$ws.Range("K136").Value = 13660.7271
$ws.Range("L136").Value = 12999
$ws.Range("M136").Value = -11110.7271
$ws.Range("N136").Value = -18099

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 32374.5
$ws.Range("I48").Value = 29750
$ws.Range("J48").Value = 34999
$ws.Range("K48").Value = 29750
$ws.Range("L48").Value = 34999
$ws.Range("M48").Value = -29181
$ws.Range("N48").Value = -36137

$ws.Range("H81").Value = 48885.453
$ws.Range("I81").Value = 68302.336
$ws.Range("J81").Value = 7277.857
$ws.Range("K81").Value = 136604.672
$ws.Range("L81").Value = 14555.714
$ws.Range("M81").Value = -135543.672
$ws.Range("N81").Value = -16677.714

$ws.Range("H84").Value = 48885.453
$ws.Range("I84").Value = 68302.336
$ws.Range("J84").Value = 7277.857
$ws.Range("K84").Value = 683023.36
$ws.Range("L84").Value = 72778.57000000001
$ws.Range("M84").Value = -677719.36
$ws.Range("N84").Value = -83386.57000000001

$ws.Range("H100").Value = 1001.65
$ws.Range("I100").Value = 884.17645
$ws.Range("J100").Value = 1667.3334
$ws.Range("K100").Value = 1768.3529
$ws.Range("L100").Value = 3334.6668
$ws.Range("M100").Value = -1227.3529
$ws.Range("N100").Value = -4416.6668

$ws.Range("H107").Value = 724.5925999999999
$ws.Range("I107").Value = 529.6842
$ws.Range("J107").Value = 1187.5
$ws.Range("K107").Value = 1589.0526
$ws.Range("L107").Value = 3562.5
$ws.Range("M107").Value = 330.9474
$ws.Range("N107").Value = -7402.5

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
